# Apply the MSGS message-table rework:
#   Row 8 ("Came for Open Day")       -> "Local - INTRO - Personalized"
#   Row 9 ("Didn't come for Open Day")-> "International - INTRO - Personalized"
#   Row 10 (old "INTRO - Personalized") is removed (its content now lives in rows 8/9)

$wb = $excel.ActiveWorkbook
$msgs = $wb.Worksheets.Item("MSGS")

$msgs.Range("B8").Value = '%5BHello+%7C+Hey+%7C+Hey+there%5D+%7Bname%7D%21+%F0%9F%91%8B%F0%9F%8F%BC%E2%9C%A8%0A%0AI%E2%80%99m+%7Bsender%7D+from+Heriot-Watt+University+Malaysia+%E2%80%93+awesome+to+connect+with+you%21+%F0%9F%9A%80%0A%0AI%27ve+received+your+interests+in+%7Bcourse%7D%2C+and+I%E2%80%99d+love+to+share+some+exciting+details+with+you%21+%F0%9F%92%A1%F0%9F%8E%93%0A%0ALet+me+know+how+you%E2%80%99d+like+to+proceed%3A%0A%0A%F0%9F%93%B5+Reply+1+%E2%80%93+Nope%2C+not+interested.+Please+stop+contacting+me.%0A%F0%9F%92%AC+Reply+2+%E2%80%93+Yes%21+Send+me+more+details.%0A%0AJust+reply+1+or+2%2C+and+I%E2%80%99ll+handle+the+rest%21+Looking+forward+to+chatting+with+you%21+%F0%9F%98%83'
$msgs.Range("C8").Value = 'Local - INTRO - Personalized'
$msgs.Range("B9").Value = '%5BHello+%7C+Hey+%7C+Hey+there%5D+%7Bname%7D%21+%F0%9F%91%8B%F0%9F%8F%BC%E2%9C%A8%0A%0AI%E2%80%99m+%7Bsender%7D+from+Heriot-Watt+University+Malaysia+%E2%80%93+awesome+to+connect+with+you%21+%F0%9F%9A%80%0A%0AI%27ve+received+your+interests+in+%7Bcourse%7D%2CIf+you+would+like+to+enquire%2C+please+find+our+designated+consultants+in+charge+for+international+student+enquiry+and+affairs+Ms.+Jane+%2B60+17-227+3699+or+Mr.+Amir+%2B60+12-953+0199+for+further+enquiry+and+they%27ll+be+able+to+help+your+further.'
$msgs.Range("C9").Value = 'International - INTRO - Personalized'

# Drop the now-redundant last row of the Message table
$msgs.Rows("10:10").Delete()

# Column C grew a bit to fit the new "Info" labels
$msgs.Columns("C:C").ColumnWidth = 29.49

$msgs.Range("D9").Select()
$msgs.Activate()

# The LIST sheet cursor had moved to G11 before MSGS became the active tab
$list = $wb.Worksheets.Item("LIST")
$list.Range("G11").Select()

$msgs.Activate()
